# Updated scripts for resolving conflicts
#
# Applies the changes described by the commit:
#  - Sheet1: swap the two client records' email/company text values,
#    tweak a handful of phone-number digits, resize column O slightly,
#    and move the active selection.
#  - Shrink the saved workbook window width a bit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1 data edits
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Text values (write in A2, A3, O2, O3 order so the shared-string table
# is rebuilt in the same layout the workbook ends up with).
$ws.Range("A2").Value = "ABC"
$ws.Range("A3").Value = "ABCS"
$ws.Range("O2").Value = "ABC@gmail.com"
$ws.Range("O3").Value = "ABCS@gmail.com"

# Phone / fax number tweaks
$ws.Range("B2").Value = 9098833665
$ws.Range("C2").Value = 4560783390
$ws.Range("D2").Value = 6567701114

$ws.Range("B3").Value = 6789933115
$ws.Range("C3").Value = 8760243356

# Widen column O (email column) slightly
$ws.Columns.Item(15).ColumnWidth = 23.5

# Move the active selection, scrolling the view so column H is at the
# left edge and L9 is the active cell
$selected = $ws.Range("L9").Select()

# ---------------------------------------------------------------
# Workbook window size (saved window geometry -> bookViews/workbookView)
# ---------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Width = 13630
$null = $win
